$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet's related name (sheet tab name) from SCD0323 to SCD0024
$ws.Name = "SCD0024"

# Update the TC_ID value in B2 from "DGS-338" to "SCD0024-002"
$ws.Range("B2").Value = "SCD0024-002"

# Column B needs to widen to fit the new (longer) content - closest achievable
# to the authored 12.7109375 character width given this host's width grid.
$ws.Columns.Item(2).ColumnWidth = 11.877604166666666

# Update selection / view state: active cell now B3, scrolled back to top-left A1
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
